# Applies the cryptos.xlsx data refresh (Wed Jun 19 13:00:18 UTC 2024 GitHub Actions run).
# Updates Price (D) and Volume(1h) (E) columns for most rows, and for a few
# rows (33-35, 38-39) the coin moved rank so Coin (B) / Link (C) / Price (D)
# were all replaced with another coin's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
  @{ Row=2; B=$null; C=$null; D="65.173.56"; E="  -0.04%  " },
  @{ Row=3; B=$null; C=$null; D="3.546.75"; E="  +4.20%  " },
  @{ Row=4; B=$null; C=$null; D=$null; E="  +0.04%  " },
  @{ Row=5; B=$null; C=$null; D="598.44"; E="  +3.11%  " },
  @{ Row=6; B=$null; C=$null; D="138.57"; E="  +1.90%  " },
  @{ Row=7; B=$null; C=$null; D="3.546.97"; E="  +4.30%  " },
  @{ Row=8; B=$null; C=$null; D=$null; E="  +0.16%  " },
  @{ Row=9; B=$null; C=$null; D="0.495"; E="  +0.32%  " },
  @{ Row=10; B=$null; C=$null; D=$null; E="  +3.78%  " },
  @{ Row=11; B=$null; C=$null; D="6.93"; E="  -3.00%  " },
  @{ Row=12; B=$null; C=$null; D="0.387"; E="  +4.12%  " },
  @{ Row=13; B=$null; C=$null; D="4.148.31"; E="  +4.48%  " },
  @{ Row=14; B=$null; C=$null; D="0.0000184"; E="  +2.99%  " },
  @{ Row=15; B=$null; C=$null; D="27.37"; E="  +5.55%  " },
  @{ Row=16; B=$null; C=$null; D="3.553.19"; E="  +4.57%  " },
  @{ Row=17; B=$null; C=$null; D=$null; E="  +1.52%  " },
  @{ Row=18; B=$null; C=$null; D="65.125.84"; E="  -0.03%  " },
  @{ Row=19; B=$null; C=$null; D="10.07"; E="  +5.59%  " },
  @{ Row=20; B=$null; C=$null; D="5.91"; E="  +1.17%  " },
  @{ Row=21; B=$null; C=$null; D="14.28"; E="  +5.62%  " },
  @{ Row=22; B=$null; C=$null; D="392.75"; E="  +3.30%  " },
  @{ Row=23; B=$null; C=$null; D="0.575"; E=$null },
  @{ Row=24; B=$null; C=$null; D="3.685.08"; E="  +4.16%  " },
  @{ Row=25; B=$null; C=$null; D="73.69"; E="  +2.57%  " },
  @{ Row=26; B=$null; C=$null; D=$null; E="  -0.14%  " },
  @{ Row=27; B=$null; C=$null; D=$null; E="  +9.92%  " },
  @{ Row=28; B=$null; C=$null; D="7.87"; E="  +11.69%  " },
  @{ Row=29; B=$null; C=$null; D=$null; E="  -0.03%  " },
  @{ Row=30; B=$null; C=$null; D=$null; E="  +3.97%  " },
  @{ Row=31; B=$null; C=$null; D="8.34"; E="  +4.21%  " },
  @{ Row=32; B=$null; C=$null; D="3.566.53"; E="  +4.55%  " },
  @{ Row=33; B="Fetch.AI"; C="https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"; D="1.38"; E="  +20.41%  " },
  @{ Row=34; B="USDe"; C="https://coinranking.com/coin/exbfr2U-0+usde-usde"; D="1.00"; E="  +0.02%  " },
  @{ Row=35; B="EthereumClassic"; C="https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; D="23.86"; E="  +4.91%  " },
  @{ Row=36; B=$null; C=$null; D="0.145"; E="  +1.96%  " },
  @{ Row=37; B=$null; C=$null; D="1.59"; E="  +8.82%  " },
  @{ Row=38; B="Aptos"; C="https://coinranking.com/coin/HGYj5JCv5+aptos-apt"; D="6.96"; E="  +4.17%  " },
  @{ Row=39; B="Monero"; C="https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; D="169.72"; E="  +0.02%  " },
  @{ Row=40; B=$null; C=$null; D="5.05"; E="  +8.27%  " },
  @{ Row=41; B=$null; C=$null; D="0.0806"; E="  +7.15%  " },
  @{ Row=42; B=$null; C=$null; D="0.823"; E="  +2.52%  " },
  @{ Row=43; B=$null; C=$null; D="26.65"; E="  +21.51%  " },
  @{ Row=44; B=$null; C=$null; D="42.51"; E="  -1.69%  " },
  @{ Row=45; B=$null; C=$null; D="0.999"; E="  +0.09%  " },
  @{ Row=46; B=$null; C=$null; D="4.45"; E="  +2.75%  " },
  @{ Row=47; B=$null; C=$null; D="1.21"; E="  +10.46%  " },
  @{ Row=48; B=$null; C=$null; D=$null; E="  +5.44%  " },
  @{ Row=49; B=$null; C=$null; D="6.85"; E="  +6.19%  " },
  @{ Row=50; B=$null; C=$null; D="2.382.54"; E="  +10.33%  " },
  @{ Row=51; B=$null; C=$null; D="309.40"; E="  +16.22%  " }
)

foreach ($item in $updates) {
    $r = $item.Row

    foreach ($col in @("B", "C", "E")) {
        $newVal = $item[$col]
        if ($newVal -ne $null) {
            $ws.Range("$col$r").Value = $newVal
        }
    }

    # Column D holds free-form price text (e.g. "3.546.75", "65.173.56") that
    # must stay a text string even when it would otherwise parse as a number
    # (e.g. "598.44", "1.38"). Force the cell to Text format before writing,
    # then restore the default "Normal" style so no stray formatting is left
    # behind on cells that originally had no explicit style.
    $newD = $item["D"]
    if ($newD -ne $null) {
        $addr = "D$r"
        $ws.Range($addr).NumberFormat = "@"
        $ws.Range($addr).Value = $newD
        $ws.Range($addr).Style = "Normal"
    }
}

